# Apply requested numeric corrections to table_12_ft table (L2C Q4 2021 report draft 2)
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "217 (72.1)" "219 (72.0)"
Replace-Text "69 (68.3)"  "69 (67.6)"
Replace-Text "72 (73.5)"  "74 (74.0)"
Replace-Text "32 (10.6)"  "32 (10.5)"
Replace-Text "11 (10.9)"  "11 (10.8)"
Replace-Text "8 (8.2)"    "8 (8.0)"
Replace-Text "143 (47.5)" "144 (47.4)"
Replace-Text "45 (44.6)"  "45 (44.1)"
Replace-Text "51 (52.0)"  "52 (52.0)"
Replace-Text "68 (22.6)"  "69 (22.7)"
Replace-Text "27 (26.7)"  "28 (27.5)"
Replace-Text "21 (21.4)"  "21 (21.0)"
